# Updated cryptos list on Sun Sep  8 18:56:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.281.20"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.266.91"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "496.21"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "128.93"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").Value = "22.91"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "2.668.91"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "54.242.70"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "2.278.45"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "10.25"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "301.76"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D23").Value = "60.96"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "7.33"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").Value = "170.99"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").Value = "1.08"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "17.81"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("E35").Value = "  +8.95%  "
$ws.Range("D36").Value = "1.19"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "125.26"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").Value = "0.0895"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "0.547"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "241.56"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "16.12"
$ws.Range("E50").Value = "  -0.59%  "

# Row 51 changed from BitgetToken to ZEEBU
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").Value = "4.59"
$ws.Range("E51").Value = "  -1.05%  "
